$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells in column D whose new values look numeric need to be forced to text
# (matching the source which stores them as inline strings), so we temporarily
# apply a text number format before assigning the values.
$textForceCells = @('D5', 'D14', 'D19', 'D27', 'D29', 'D31', 'D32', 'D36', 'D38', 'D40', 'D44', 'D45', 'D46', 'D48')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped from the refreshed crypto feed
$ws.Range('D2').Value = '63.710.04'
$ws.Range('E2').Value = '  -3.24%  '
$ws.Range('D3').Value = '2.613.23'
$ws.Range('E3').Value = '  -1.74%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '574.91'
$ws.Range('E5').Value = '  -4.13%  '
$ws.Range('E6').Value = '  -2.14%  '
$ws.Range('E8').Value = '  -2.90%  '
$ws.Range('D9').Value = '2.611.19'
$ws.Range('E9').Value = '  -1.70%  '
$ws.Range('E10').Value = '  -6.75%  '
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('E12').Value = '  -5.12%  '
$ws.Range('E13').Value = '  -0.37%  '
$ws.Range('D14').Value = '28.23'
$ws.Range('E14').Value = '  -3.05%  '
$ws.Range('D15').Value = '3.085.84'
$ws.Range('E15').Value = '  -1.79%  '
$ws.Range('E16').Value = '  -8.07%  '
$ws.Range('D17').Value = '63.624.74'
$ws.Range('E17').Value = '  -3.19%  '
$ws.Range('D18').Value = '2.626.23'
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('D19').Value = '12.00'
$ws.Range('E19').Value = '  -4.64%  '
$ws.Range('E20').Value = '  +1.77%  '
$ws.Range('E21').Value = '  -5.75%  '
$ws.Range('E22').Value = '  -2.94%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('E24').Value = '  -3.63%  '
$ws.Range('E25').Value = '  +0.32%  '
$ws.Range('E26').Value = '  -3.54%  '
$ws.Range('D27').Value = '592.52'
$ws.Range('E27').Value = '  +2.00%  '
$ws.Range('E28').Value = '  -4.84%  '
$ws.Range('D29').Value = '1.57'
$ws.Range('E29').Value = '  -2.90%  '
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').Value = '7.90'
$ws.Range('E32').Value = '  -2.92%  '
$ws.Range('E33').Value = '  -4.24%  '
$ws.Range('E34').Value = '  -4.12%  '
$ws.Range('E35').Value = '  -2.60%  '
$ws.Range('D36').Value = '5.37'
$ws.Range('E36').Value = '  -3.00%  '
$ws.Range('E37').Value = '  -4.73%  '
$ws.Range('D38').Value = '19.73'
$ws.Range('E38').Value = '  -4.30%  '
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('D40').Value = '154.15'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('E41').Value = '  -4.47%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('E43').Value = '  +0.70%  '
$ws.Range('D44').Value = '41.42'
$ws.Range('E44').Value = '  -3.43%  '
$ws.Range('D45').Value = '157.55'
$ws.Range('E45').Value = '  -2.58%  '
$ws.Range('D46').Value = '23.78'
$ws.Range('E46').Value = '  +1.76%  '
$ws.Range('E47').Value = '  -5.12%  '
$ws.Range('D48').Value = '0.0589'
$ws.Range('E48').Value = '  -4.53%  '
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('E50').Value = '  -1.24%  '
$ws.Range('E51').Value = '  -5.20%  '

# Restore the default (Normal) style on the cells we forced to text so no
# stray formatting is introduced
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
